$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 256.85715
$ws.Range("I12").Value = 129.2
$ws.Range("J12").Value = 576
$ws.Range("K12").Value = 129.2
$ws.Range("L12").Value = 576
$ws.Range("M12").Value = 40.80000000000001
$ws.Range("N12").Value = -916

$ws.Range("H33").Value = 165.4
$ws.Range("I33").Value = 165.4
$ws.Range("K33").Value = 165.4
$ws.Range("M33").Value = 63.59999999999999

$ws.Range("H64").Value = 4033.25

$ws.Range("H67").Value = 4033.25

$ws.Range("H138").Value = 1927.7222
$ws.Range("J138").Value = 3714.2856
$ws.Range("L138").Value = 11142.8568
$ws.Range("N138").Value = -21422.8568

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 3099.75
$ws.Range("I2").Value = 1466.3334
$ws.Range("J2").Value = 8000
$ws.Range("K2").Value = 1466.3334
$ws.Range("L2").Value = 8000
$ws.Range("M2").Value = -1353.3334
$ws.Range("N2").Value = -8226

$ws.Range("H45").Value = 2166.6667
$ws.Range("J45").Value = 2500
$ws.Range("L45").Value = 2500
$ws.Range("N45").Value = -3254

$ws.Range("H116").Value = 3099.75
$ws.Range("I116").Value = 1466.3334
$ws.Range("J116").Value = 8000
$ws.Range("K116").Value = 1466.3334
$ws.Range("L116").Value = 8000
$ws.Range("M116").Value = 827.6666
$ws.Range("N116").Value = -12588

$ws.Range("H122").Value = 0
$ws.Range("I122").Value = 0
$ws.Range("K122").Value = 0
$ws.Range("M122").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 3099.75
$ws.Range("I3").Value = 1466.3334
$ws.Range("J3").Value = 8000
$ws.Range("K3").Value = 1466.3334
$ws.Range("L3").Value = 8000
$ws.Range("M3").Value = -1352.3334
$ws.Range("N3").Value = -8228

$ws.Range("H20").Value = 757.6
$ws.Range("I20").Value = 696.3333
$ws.Range("K20").Value = 696.3333
$ws.Range("M20").Value = -449.3333

$ws.Range("H105").Value = 1495
$ws.Range("J105").Value = 0
$ws.Range("L105").Value = 0
$ws.Range("N105").ClearContents()

$ws.Range("H134").Value = 7530.2085
$ws.Range("J134").Value = 16500
$ws.Range("L134").Value = 49500
$ws.Range("N134").Value = -54570

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H3").Value = 100
$ws.Range("I3").Value = 0
$ws.Range("J3").Value = 100
$ws.Range("K3").Value = 0
$ws.Range("L3").Value = 100
$ws.Range("M3").ClearContents()
$ws.Range("N3").Value = -326

$ws.Range("H31").Value = 3618.125
$ws.Range("I31").Value = 2749.5652
$ws.Range("K31").Value = 2749.5652
$ws.Range("M31").Value = -2454.5652

$ws.Range("H34").Value = 3618.125
$ws.Range("I34").Value = 2749.5652
$ws.Range("K34").Value = 2749.5652
$ws.Range("M34").Value = -2547.5652

$ws.Range("H48").Value = 12750
$ws.Range("J48").Value = 12750
$ws.Range("L48").Value = 12750
$ws.Range("N48").Value = -13702

$ws.Range("H134").Value = 3511.4285
$ws.Range("I134").Value = 3538.5454
$ws.Range("K134").Value = 10615.6362
$ws.Range("M134").Value = -8080.636200000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H36").Value = 443
$ws.Range("I36").Value = 443
$ws.Range("J36").Value = 0
$ws.Range("K36").Value = 1329
$ws.Range("L36").Value = 0
$ws.Range("M36").Value = -1160
$ws.Range("N36").ClearContents()

$ws.Range("H109").Value = 911.7778
$ws.Range("I109").Value = 911.7778
$ws.Range("J109").Value = 0
$ws.Range("K109").Value = 2735.3334
$ws.Range("L109").Value = 0
$ws.Range("M109").Value = -1695.3334
$ws.Range("N109").ClearContents()

$ws.Range("H140").Value = 1999.6666
$ws.Range("I140").Value = 1449.6
$ws.Range("K140").Value = 4348.799999999999
$ws.Range("M140").Value = 831.2000000000007

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 5848.2
$ws.Range("I113").Value = 4500
$ws.Range("K113").Value = 4500
$ws.Range("M113").Value = -2330

$ws.Range("H122").Value = 12502821
$ws.Range("I122").Value = 15626775
$ws.Range("K122").Value = 46880325
$ws.Range("M122").Value = -46877875

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1125
$ws.Range("I22").Value = 1125
$ws.Range("K22").Value = 1125
$ws.Range("M22").Value = -830

$ws.Range("H27").Value = 1125
$ws.Range("I27").Value = 1125
$ws.Range("K27").Value = 1125
$ws.Range("M27").Value = -1018

$ws.Range("H40").Value = 2460.3333
$ws.Range("I40").Value = 2359.2
$ws.Range("K40").Value = 2359.2
$ws.Range("M40").Value = -2223.2

$ws.Range("H45").Value = 12997
$ws.Range("J45").Value = 12997
$ws.Range("L45").Value = 12997
$ws.Range("N45").Value = -13811

$ws.Range("H53").Value = 19887.5
$ws.Range("I53").Value = 19887.5
$ws.Range("J53").Value = 0
$ws.Range("K53").Value = 19887.5
$ws.Range("L53").Value = 0
$ws.Range("M53").Value = -19369.5
$ws.Range("N53").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H3").Value = 0
$ws.Range("I3").Value = 0
$ws.Range("K3").Value = 0
$ws.Range("M3").ClearContents()

$ws.Range("H43").Value = 25000
$ws.Range("J43").Value = 25000
$ws.Range("L43").Value = 25000
$ws.Range("N43").Value = -25298

$ws.Range("H50").Value = 30000
$ws.Range("J50").Value = 30000
$ws.Range("L50").Value = 30000
$ws.Range("N50").Value = -31262

$ws.Range("H52").Value = 0
$ws.Range("I52").Value = 0
$ws.Range("K52").Value = 0
$ws.Range("M52").ClearContents()

$ws.Range("H53").Value = 0
$ws.Range("I53").Value = 0
$ws.Range("K53").Value = 0
$ws.Range("M53").ClearContents()

$ws.Range("H54").Value = 30000
$ws.Range("J54").Value = 30000
$ws.Range("L54").Value = 30000
$ws.Range("N54").Value = -31040

$ws.Range("H81").Value = 0
$ws.Range("I81").Value = 0
$ws.Range("J81").Value = 0
$ws.Range("K81").Value = 0
$ws.Range("L81").Value = 0
$ws.Range("M81").ClearContents()
$ws.Range("N81").ClearContents()

$ws.Range("H84").Value = 0
$ws.Range("I84").Value = 0
$ws.Range("J84").Value = 0
$ws.Range("K84").Value = 10000
$ws.Range("L84").Value = 4500
$ws.Range("M84").ClearContents()
$ws.Range("N84").ClearContents()

$ws.Range("H107").Value = 498.875
$ws.Range("I107").Value = 498.7143
$ws.Range("K107").Value = 1496.1429
$ws.Range("M107").Value = 423.8571000000002

$ws.Range("H132").Value = 1766.3572
$ws.Range("I132").Value = 1673.091
$ws.Range("K132").Value = 5019.272999999999
$ws.Range("M132").Value = -2489.272999999999
